# NIT-9008596341 "Estado de Cuenta": update the workers' arrears period
# (Periodo Mora) and basic salary (Salario Basico) figures for the new
# statement cycle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodo Mora (column E, rows 16-18)
$ws.Range("E16").Value = "2411"
$ws.Range("E17").Value = "2412"
$ws.Range("E18").Value = "2501"

# Salario Basico (column G, rows 16-18)
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
$ws.Range("G18").Value = 1423500
